# Commit: "Code Merge Changes- Xls files are updated - 9/25/2017"
#
# 1. Update the shared TestResultExcelFilePath value (cell H2) used by the
#    several "ProcessPayrollFor...SAPP" / report sheets - rename the result
#    file from "Automation Test Result for Statutory Scenarios201718.xlsx"
#    to "201718 Automation TestResult For Statutory Scenarios.xlsx".
# 2. Move the remembered cell-selection on three of the report sheets.

$wb = $excel.ActiveWorkbook

$newPath = "F:\\Automation_TestResults\\Payroll_Tax_StatutoryScenarios\\201718 Automation TestResult For Statutory Scenarios.xlsx"

$sheetsWithPath = @(
    "ProcessPayrollForJulyMonthSAPP",
    "ProcessPayrollForAugMonthSAPP",
    "ProcessPayrollForSepMonthSAPP",
    "AverageWeeklyEarningsTestReport",
    "ProcessPayrollForJan16MonthSAPP"
)

foreach ($name in $sheetsWithPath) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = $newPath
}

# Update the remembered selections (activeCell) on the affected sheets.
$ws10 = $wb.Worksheets.Item("AverageWeeklyEarningsTestReport")
$ws10.Activate()
$ws10.Range("H2").Select()

$ws7 = $wb.Worksheets.Item("ProcessPayrollForAugMonthSAPP")
$ws7.Activate()
$ws7.Range("H2").Select()

$ws8 = $wb.Worksheets.Item("ProcessPayrollForSepMonthSAPP")
$ws8.Activate()
$ws8.Range("H6").Select()

# Restore the originally active sheet/tab so we don't disturb the workbook's
# remembered active tab as a side effect of the selection changes above.
$ws11 = $wb.Worksheets.Item("ProcessPayrollForJan16MonthSAPP")
$ws11.Activate()
